$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "FUNCIONÁRIO 1"
$ws.Range("B3").Value = "FUNCIONÁRIO 2"
$ws.Range("B4").Value = "FUNCIONÁRIO 3"
$ws.Range("B5").Value = "FUNCIONÁRIO 4"
$ws.Range("B6").Value = "FUNCIONÁRIO 5"

$ws.Range("E10").Select()
